$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPS Data")

# --- Update target values (P and R columns) for rows 2-10 ---
$ws.Range("P2").Value = 60
$ws.Range("R2").Value = 100

$ws.Range("P3").Value = 30

$ws.Range("P4").Value = 45
$ws.Range("R4").Value = 68

$ws.Range("P5").Value = 30
$ws.Range("R5").Value = 56

$ws.Range("P6").Value = 45
$ws.Range("R6").Value = 68

$ws.Range("P7").Value = 30
$ws.Range("R7").Value = 56

$ws.Range("P8").Value = 40
$ws.Range("R8").Value = 60

$ws.Range("P9").Value = 40

$ws.Range("P10").Value = 38
$ws.Range("R10").Value = 44

# --- Update selection to match the new active view ---
$ws.Range("P1:S10").Select()
